$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a plain decimal number (single "." , digits only)
# must be forced to Text format first so Excel does not silently convert them
# into a floating point number (which would also round-trip imprecisely and
# drop significant digits such as trailing zeros).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "27.289.84"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").Value = "1.908.49"
$ws.Range("E3").Value = "  +2.22%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "307.92"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("E6").Value = "  +0.04%  "

Set-TextValue $ws.Range("D7") "0.5253"
$ws.Range("E7").Value = "  +3.25%  "

Set-TextValue $ws.Range("D8") "0.3788"
$ws.Range("E8").Value = "  +3.57%  "

Set-TextValue $ws.Range("D9") "0.07265"
$ws.Range("E9").Value = "  +1.25%  "

Set-TextValue $ws.Range("D10") "21.35"
$ws.Range("E10").Value = "  +3.83%  "

Set-TextValue $ws.Range("D11") "0.9020"

Set-TextValue $ws.Range("D12") "0.08291"
$ws.Range("E12").Value = "  +10.72%  "

$ws.Range("D13").Value = "1.911.02"
$ws.Range("E13").Value = "  +2.22%  "

Set-TextValue $ws.Range("D14") "95.21"
$ws.Range("E14").Value = "  +0.79%  "

Set-TextValue $ws.Range("D15") "5.296"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("E16").Value = "  -0.06%  "

Set-TextValue $ws.Range("D17") "0.000008611"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("E18").Value = "  +2.47%  "

Set-TextValue $ws.Range("D19") "1.000"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").Value = "27.323.26"
$ws.Range("E20").Value = "  +1.43%  "

Set-TextValue $ws.Range("D21") "5.075"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("D22").Value = "2.154.33"

$ws.Range("E23").Value = "  +3.16%  "

Set-TextValue $ws.Range("D24") "6.467"
$ws.Range("E24").Value = "  +1.57%  "

Set-TextValue $ws.Range("D25") "2.306"
$ws.Range("E25").Value = "  +10.22%  "

Set-TextValue $ws.Range("D26") "146.17"
$ws.Range("E26").Value = "  -0.85%  "

$ws.Range("E27").Value = "  -1.74%  "

Set-TextValue $ws.Range("D28") "18.19"
$ws.Range("E28").Value = "  +1.88%  "

Set-TextValue $ws.Range("D29") "115.11"
$ws.Range("E29").Value = "  +1.37%  "

Set-TextValue $ws.Range("D30") "5.001"
$ws.Range("E30").Value = "  +6.50%  "

Set-TextValue $ws.Range("D31") "4.816"
$ws.Range("E31").Value = "  +2.75%  "

$ws.Range("E32").Value = "  +1.44%  "

Set-TextValue $ws.Range("D33") "0.8076"
$ws.Range("E33").Value = "  +7.90%  "

Set-TextValue $ws.Range("D35") "1.244"
$ws.Range("E35").Value = "  +7.94%  "

Set-TextValue $ws.Range("D36") "2.976"
$ws.Range("E36").Value = "  +0.36%  "

Set-TextValue $ws.Range("D37") "3.365"
$ws.Range("E37").Value = "  +4.68%  "

Set-TextValue $ws.Range("D38") "2.574"
$ws.Range("E38").Value = "  +2.81%  "

Set-TextValue $ws.Range("D39") "0.5728"
$ws.Range("E39").Value = "  +3.26%  "

Set-TextValue $ws.Range("D40") "0.01985"
$ws.Range("E40").Value = "  -0.10%  "

Set-TextValue $ws.Range("D41") "1.075"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D42") "9.021"
$ws.Range("E42").Value = "  +5.15%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D43") "119.74"
$ws.Range("E43").Value = "  +3.39%  "

$ws.Range("E44").Value = "  +0.79%  "

Set-TextValue $ws.Range("D45") "0.1518"
$ws.Range("E45").Value = "  +2.35%  "

Set-TextValue $ws.Range("D46") "0.4839"
$ws.Range("E46").Value = "  +1.96%  "

Set-TextValue $ws.Range("D47") "10.17"
$ws.Range("E47").Value = "  +0.75%  "

Set-TextValue $ws.Range("D48") "1.001"
$ws.Range("E48").Value = "  +0.06%  "

Set-TextValue $ws.Range("D49") "1.619"
$ws.Range("E49").Value = "  +4.11%  "

Set-TextValue $ws.Range("D50") "37.64"
$ws.Range("E50").Value = "  +1.45%  "

Set-TextValue $ws.Range("D51") "63.83"
$ws.Range("E51").Value = "  +1.35%  "
